$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 15 (duplicate Rifampicin 5mg/kg row) - this was causing a "double box"
# for the same dosage/drug. Deleting the row shifts rows 16-23 up to 15-22.
$ws.Rows.Item(15).Delete()

# Reflect the resulting selection (the whole row that is now row 15) and reset
# the scroll position of the sheet view.
$ws.Rows.Item(15).Select() | Out-Null
